# Apply the update described by the diff:
#  - Old row 39 (Artomyces cristatus / Tommy Solberg) shifts down to row 41.
#  - Two brand-new observation rows are inserted at 39-40 (Anton Larsson).
#  - One brand-new observation row is appended at the end, row 42 (Anton Larsson).
#  - Sheet dimension grows from A1:AY39 to A1:AY42 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert two blank rows above the current row 39. This pushes
#    the existing row 39 (and everything below it) down to row 41, keeping
#    all of its cell values/types intact.
# ---------------------------------------------------------------------------
$ws.Range("A39:A40").EntireRow.Insert()

# ---------------------------------------------------------------------------
# Helper to write one observation record into a given row without Excel's
# autodetection turning the Startdatum/Slutdatum text into real dates.
# Parameters are positional (named binding is not supported by this host):
#   1 Row, 2 Id, 3 TaxonSort, 4 Status, 5 RedList, 6 TaxonId, 7 ArtNamn,
#   8 VetNamn, 9 Auktor, 10 Lokalnamn, 11 Ost, 12 Nord, 13 Noggrannhet,
#   14 Lan, 15 Kommun, 16 Provins, 17 Forsamling, 18 StartDatum, 19 SlutDatum,
#   20 Substratnamn, 21 VetSubstratnamn, 22 SubstratBeskrivning,
#   23 Rapportor, 24 Observatorer
# ---------------------------------------------------------------------------
function Set-ArtRow($Row, $Id, $TaxonSort, $Status, $RedList, $TaxonId, $ArtNamn, $VetNamn, $Auktor, $Lokalnamn, $Ost, $Nord, $Noggrannhet, $Lan, $Kommun, $Provins, $Forsamling, $StartDatum, $SlutDatum, $Substratnamn, $VetSubstratnamn, $SubstratBeskrivning, $Rapportor, $Observatorer) {

    $ws.Cells.Item($Row, 1).Value = $Id
    $ws.Cells.Item($Row, 2).Value = $TaxonSort
    $ws.Cells.Item($Row, 3).Value = $Status
    $ws.Cells.Item($Row, 4).Value = $RedList
    $ws.Cells.Item($Row, 5).Value = $TaxonId
    $ws.Cells.Item($Row, 6).Value = $ArtNamn
    $ws.Cells.Item($Row, 7).Value = $VetNamn
    $ws.Cells.Item($Row, 8).Value = $Auktor

    $ws.Cells.Item($Row, 16).Value = $Lokalnamn
    $ws.Cells.Item($Row, 17).Value = $Ost
    $ws.Cells.Item($Row, 18).Value = $Nord
    $ws.Cells.Item($Row, 19).Value = $Noggrannhet
    $ws.Cells.Item($Row, 20).Value = $Lan
    $ws.Cells.Item($Row, 21).Value = $Kommun
    $ws.Cells.Item($Row, 22).Value = $Provins
    $ws.Cells.Item($Row, 23).Value = $Forsamling

    # Columns Y (25) and AA (27) hold the dates as plain text, matching the
    # source export -- force text formatting first so "2023-09-26" isn't
    # reinterpreted as a date serial number.
    $dateCell = $ws.Cells.Item($Row, 25)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $StartDatum

    $dateCell2 = $ws.Cells.Item($Row, 27)
    $dateCell2.NumberFormat = "@"
    $dateCell2.Value = $SlutDatum

    $ws.Cells.Item($Row, 30).Value = $false   # AD - Ej aterfunnen
    $ws.Cells.Item($Row, 31).Value = $false   # AE - Osaker artbestamning
    $ws.Cells.Item($Row, 33).Value = $false   # AG - Ospontan

    if ($Substratnamn) {
        $ws.Cells.Item($Row, 36).Value = $Substratnamn        # AJ
    }
    if ($VetSubstratnamn) {
        $ws.Cells.Item($Row, 37).Value = $VetSubstratnamn     # AK
    }
    if ($SubstratBeskrivning) {
        $ws.Cells.Item($Row, 41).Value = $SubstratBeskrivning # AO
    }

    $ws.Cells.Item($Row, 49).Value = $Rapportor       # AW
    $ws.Cells.Item($Row, 50).Value = $Observatorer    # AX
}

# ---------------------------------------------------------------------------
# 2) Row 39 - new record: Bazzania trilobata (Stor revmossa)
# ---------------------------------------------------------------------------
Set-ArtRow 39 112387478 95006 "Ovaliderad" "LC" 2569 "Stor revmossa" "Bazzania trilobata" "(L.) Gray" "Stora Stickshöjden, Dls" 318275 6554953 5 "Västra Götaland" "Dals-Ed" "Dalsland" "Nössemark" "2023-09-26" "2023-09-26" $null $null $null "Anton Larsson" "Anton Larsson, Maria Johansson"

# ---------------------------------------------------------------------------
# 3) Row 40 - new record: Odontoschisma denudatum (Kornknutmossa)
# ---------------------------------------------------------------------------
Set-ArtRow 40 112387479 94326 "Ovaliderad" "LC" 2590 "Kornknutmossa" "Odontoschisma denudatum" "(Mart.) Dumort" "Stora Stickshöjden, Dls" 318301 6554977 5 "Västra Götaland" "Dals-Ed" "Dalsland" "Nössemark" "2023-09-26" "2023-09-26" "tall" "Pinus sylvestris" "Pinus sylvestris" "Anton Larsson" "Anton Larsson, Maria Johansson"

# ---------------------------------------------------------------------------
# 4) Row 41 already contains the shifted-down original record untouched.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5) Row 42 - new record appended at the end: Hydnellum suaveolens (Dofttaggsvamp)
# ---------------------------------------------------------------------------
Set-ArtRow 42 112387492 90818 "Ovaliderad" "NT" 4368 "Dofttaggsvamp" "Hydnellum suaveolens" "(Scop.:Fr.) P. Karst." "Stora Stickshöjden, Dls" 318046 6554741 5 "Västra Götaland" "Dals-Ed" "Dalsland" "Nössemark" "2023-09-26" "2023-09-26" $null $null $null "Anton Larsson" "Anton Larsson, Maria Johansson"
